$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Shared-string bookkeeping: rename L2's current text ("Leads Prospek
#    NULL 03", the sole reference to that shared string) to "QUERY1" first so
#    the existing slot is reused in place; T1 will then reuse this same
#    "QUERY1" string. L2 gets repointed to a brand-new "Leads Prospek NULL 09"
#    string at the very end of this script, once all the other new strings
#    have been appended, so it lands as the final shared-string entry.
# ---------------------------------------------------------------------------
$ws.Range("L2").Value = "QUERY1"

# ---------------------------------------------------------------------------
# 2) New header row (row 1) cells T1:AB1
# ---------------------------------------------------------------------------
$ws.Range("T1").Value = "QUERY1"
$ws.Range("U1").Value = "QUERY2"
$ws.Range("V1").Value = "QUERY3"
$ws.Range("W1").Value = "EXPL_QUERY1"
$ws.Range("X1").Value = "EXPL_QUERY2"
$ws.Range("Y1").Value = "EXPL_QUERY3"
$ws.Range("Z1").Value = "USER_DB"
$ws.Range("AA1").Value = "PASSWORD_DB"
$ws.Range("AB1").Value = "HOSTNAME"

# ---------------------------------------------------------------------------
# 3) New data row (row 2) formatting + values
#    T2:V2  -> wrap-text only style, T2 holds the SQL query text
#    W2:Y2  -> text-format + vertical-center + wrap style, W2 holds the note
#    Z2     -> copy G2's format (Arial font) then right-align -> "sa"
#    AA2    -> reuse G2's exact existing style (font+wrap+vcenter) -> password
#    AB2    -> reuse K2's exact existing style (vcenter only) -> hostname
# ---------------------------------------------------------------------------
$ws.Range("T2:V2").WrapText = $true
$ws.Range("T2").Value = "SELECT DISTINCT b.Npp, c.KodeOutlet, c.Name, d.KodeOutlet AS KODE_OUTLET_BNI_MULTIFINANCE FROM DigisalesNew..Tbl_Pegawai AS b join DigisalesNew..Tbl_Unit AS c ON b.Unit_Id = c.Id Left Join Digisales_Leads..MappingBniMultifinance as d ON c.KodeOutlet = d.KodeOutlet WHERE Npp = '22914' OR Npp = '49998'"

$ws.Range("W2:Y2").NumberFormat = "@"
$ws.Range("W2:Y2").VerticalAlignment = -4108
$ws.Range("W2:Y2").WrapText = $true

$ws.Range("G2").Copy()
$ws.Range("Z2").PasteSpecial(-4122)
$ws.Range("Z2").HorizontalAlignment = -4152
$ws.Range("Z2").Value = "sa"

$ws.Range("G2").Copy()
$ws.Range("AA2").PasteSpecial(-4122)
$ws.Range("AA2").Value = "4eFfEJAA!"

$ws.Range("K2").Copy()
$ws.Range("AB2").PasteSpecial(-4122)
$ws.Range("AB2").Value = "192.168.232.6"

$ws.Range("W2").Value = "Bukti bahwa npp 49998 tidak memiliki kode outlet"

# ---------------------------------------------------------------------------
# 4) Re-point L2 to the brand-new, final shared string (must happen after all
#    of the above so it is appended last in the shared-strings table).
# ---------------------------------------------------------------------------
$ws.Range("L2").Value = "Leads Prospek NULL 09"

# ---------------------------------------------------------------------------
# 5) Column widths for the newly-used columns (best achievable given the
#    engine's column-width rounding granularity).
# ---------------------------------------------------------------------------
$ws.Columns.Item(20).ColumnWidth = 89.5
$ws.Columns.Item(23).ColumnWidth = 14.666666666666666
$ws.Range("X1:Y1").ColumnWidth = 12.333333333333332
$ws.Columns.Item(27).ColumnWidth = 13.833333333333332
$ws.Columns.Item(28).ColumnWidth = 11.833333333333332

# ---------------------------------------------------------------------------
# 6) Selection / scroll position updates.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 12
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("L2").Select()
